# Commit: "rename 27001 and 27002 tab"
# The sheet tab literally named "ISO27001+27002" is renamed to "27001+27002".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ISO27001+27002")
$ws.Name = "27001+27002"

# Keep/restore this sheet as the active one, scrolled & selected the way the
# author last left it (topLeftCell A166, selection C185) before saving.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 166
$ws.Range("C185").Select()
